# Insert a new weekly price record as row 59 in the "Poroto verde" sheet.
# This pushes the existing rows 59-104 down to become rows 60-105
# (their content stays exactly the same), and fills the newly inserted
# row 59 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 59, shifting rows 59:104 -> 60:105
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with the new record
$ws.Cells.Item(59, 1).Value  = 7
$ws.Cells.Item(59, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value  = "Ñuble"
$ws.Cells.Item(59, 4).Value  = 44907
$ws.Cells.Item(59, 5).Value  = 16
$ws.Cells.Item(59, 6).Value  = 100112031
$ws.Cells.Item(59, 7).Value  = "Poroto verde"
$ws.Cells.Item(59, 8).Value  = "Sin especificar"
$ws.Cells.Item(59, 9).Value  = "Primera"
$ws.Cells.Item(59, 10).Value = 50
$ws.Cells.Item(59, 11).Value = 35000
$ws.Cells.Item(59, 12).Value = 35000
$ws.Cells.Item(59, 13).Value = 35000
$ws.Cells.Item(59, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(59, 15).Value = "Región del Maule"
$ws.Cells.Item(59, 16).Value = 1400
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
